$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values need to be swapped between the paired rows.
$cols = @("A","B","E","F","G","H","Q","R","AJ","AK","AO")

# Row pairs that were swapped in the source data.
$pairs = @(
    @(16, 17),
    @(19, 20),
    @(26, 27)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
